# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation".
#    The same shared string is used by Overview!E2, Overview!F2,
#    zh-cn!C2 and de-de!C2, so all four cells are updated so the
#    workbook again stores a single, shared piece of text.
# 2) Narrow the "Status" columns (Overview columns E/F, and column C on
#    the zh-cn / de-de detail sheets) to better fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Update status values -------------------------------------------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- Resize the Status columns --------------------------------------------
# Target stored width is 13.4101845877511; Excel's ColumnWidth setter snaps
# to whole-pixel increments, so 12.5 is the input that lands on the closest
# representable width.
$overview.Range("E1:F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
